$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B7").Value = 120
$ws.Range("B8").Value = 300
$ws.Range("B9").Value = 230
$ws.Range("B10").Value = 250
$ws.Range("B11").Value = 70
$ws.Range("B12").Value = 80
$ws.Range("B13").Value = 300

$ws.Range("B13").Select()
